$wb = $excel.ActiveWorkbook

# --- Phase1 sheet: decrease probability of fetal death before 4 weeks ---
$ws1 = $wb.Worksheets.Item("Phase1")
$ws1.Range("B2").Value = 0.1
$ws1.Range("B3").Value = 0.1
$ws1.Range("B4").Value = 0.05
$ws1.Range("B5").Value = 0.05

# --- Phase2 sheet: same probability decrease ---
$ws2 = $wb.Worksheets.Item("Phase2")
$ws2.Range("C2").Value = 0.1
$ws2.Range("C3").Value = 0.1
$ws2.Range("C4").Value = 0.05
$ws2.Range("C5").Value = 0.05

# Update the selection shown on the Phase2 sheet
$ws2.Activate()
$ws2.Range("C2:C5").Select()

# --- Phase4 sheet: it is no longer the active/selected tab ---
$ws4 = $wb.Worksheets.Item("Phase4")
$ws4.Activate()
$ws4.Range("C6").Select()

# --- Phase1 becomes the active tab, with B2:B5 selected ---
$ws1.Activate()
$ws1.Range("B2:B5").Select()
